# Generate Report for Archive
#
# Reorders the three "in-flight" localization rows (6b9ff258, 965104be,
# 253d137b) across the Overview / zh-cn / de-de sheets so that the two
# files still being translated (6b9ff258, 965104be) move above the one
# that is ready for handoff (253d137b), and flips the in-progress files'
# status from "Ready for handoff" to "In Translation".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A3").Value() = "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md"
$ov.Range("B3").Value() = "In Translation"
$ov.Range("C3").Value() = "In Translation"
$ov.Range("D3").Value() = "2016-25-11 16:25:17"

$ov.Range("A4").Value() = "965104be-0ef5-4edb-82b1-facb0f37c968.md"
$ov.Range("B4").Value() = "In Translation"
$ov.Range("C4").Value() = "In Translation"
$ov.Range("D4").Value() = "2016-25-11 16:25:17"

$ov.Range("A5").Value() = "253d137b-9592-410f-9fca-d89327456d1f.md"
$ov.Range("B5").Value() = "Ready for handoff"
$ov.Range("C5").Value() = "Ready for handoff"
$ov.Range("D5").Value() = "2016-24-11 16:24:14"

# Hyperlinks: the relationship targets (rId -> URL) are untouched by this
# edit, only which cell shows which display text changes. Recreate the
# hyperlinks in the *original* URL order (ef4783de, 253d137b, 6b9ff258,
# 965104be) but anchor them on the *new* row layout, so the rIds line up
# exactly like the source workbook does.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/99cf95aca13517e6b2d75943be5e9e0527602912/e2e/ef4783de-cfb1-4ba8-ade8-09b399ca2752.md", "", "", "ef4783de-cfb1-4ba8-ade8-09b399ca2752.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4daf793447831388b2108a25df80716aaad753b6/e2e/253d137b-9592-410f-9fca-d89327456d1f.md", "", "", "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/684e0c5928f94a0012b5261340bed1e0c0d7575f/e2e/6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md", "", "", "965104be-0ef5-4edb-82b1-facb0f37c968.md")
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/684e0c5928f94a0012b5261340bed1e0c0d7575f/e2e/965104be-0ef5-4edb-82b1-facb0f37c968.md", "", "", "253d137b-9592-410f-9fca-d89327456d1f.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A3").Value() = "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md"
$zh.Range("B3").Value() = ".md"
$zh.Range("C3").Value() = "In Translation"
$zh.Range("D3").Value() = "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.0abd8208a0eb44c32d9f52395849381cb7361d5f.zh-cn.xlf"
$zh.Range("E3").Value() = "2016-03-11 16:25:14"
$zh.Range("H3").Value() = "0001-01-01 00:00:00"
$zh.Range("I3").Value() = "Include"

$zh.Range("A4").Value() = "965104be-0ef5-4edb-82b1-facb0f37c968.md"
$zh.Range("B4").Value() = ".md"
$zh.Range("C4").Value() = "In Translation"
$zh.Range("D4").Value() = "965104be-0ef5-4edb-82b1-facb0f37c968.a05aeaf165ae5501d28d8aece37eedefb60075ad.zh-cn.xlf"
$zh.Range("E4").Value() = "2016-03-11 16:25:14"
$zh.Range("H4").Value() = "0001-01-01 00:00:00"
$zh.Range("I4").Value() = "Include"

$zh.Range("A5").Value() = "253d137b-9592-410f-9fca-d89327456d1f.md"
$zh.Range("B5").Value() = ".md"
$zh.Range("C5").Value() = "Ready for handoff"
$zh.Range("D5").Value() = "253d137b-9592-410f-9fca-d89327456d1f.fc74a366d75ceb48353102f324d1f380db212c73.zh-cn.xlf"
$zh.Range("E5").Value() = "2016-03-11 16:24:11"
$zh.Range("H5").Value() = "0001-01-01 00:00:00"
$zh.Range("I5").Value() = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/99cf95aca13517e6b2d75943be5e9e0527602912/e2e/ef4783de-cfb1-4ba8-ade8-09b399ca2752.md", "", "", "ef4783de-cfb1-4ba8-ade8-09b399ca2752.md")
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/99cf95aca13517e6b2d75943be5e9e0527602912/e2e/ef4783de-cfb1-4ba8-ade8-09b399ca2752.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/60e7ad4a6dc8f9ee043f692f172a26f564ece132/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ef4783de-cfb1-4ba8-ade8-09b399ca2752.69d912e4e3936c4a3dc4afbd5ae619737a6dcfd6.zh-cn.xlf", "", "", "ef4783de-cfb1-4ba8-ade8-09b399ca2752.69d912e4e3936c4a3dc4afbd5ae619737a6dcfd6.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/79f11e1f6d2a89431fef1e00a29c4eda38b72bb0/e2e/ef4783de-cfb1-4ba8-ade8-09b399ca2752.md", "", "", "ef4783de-cfb1-4ba8-ade8-09b399ca2752.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8a399074f42a8891c9cd9afa61c605f124ea98f4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ef4783de-cfb1-4ba8-ade8-09b399ca2752.69d912e4e3936c4a3dc4afbd5ae619737a6dcfd6.zh-cn.xlf", "", "", "ef4783de-cfb1-4ba8-ade8-09b399ca2752.69d912e4e3936c4a3dc4afbd5ae619737a6dcfd6.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4daf793447831388b2108a25df80716aaad753b6/e2e/253d137b-9592-410f-9fca-d89327456d1f.md", "", "", "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md")
$zh.Hyperlinks.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/4daf793447831388b2108a25df80716aaad753b6/e2e/253d137b-9592-410f-9fca-d89327456d1f.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dabd5264373e6c580524dcffc1c65dc479913ba3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/253d137b-9592-410f-9fca-d89327456d1f.fc74a366d75ceb48353102f324d1f380db212c73.zh-cn.xlf", "", "", "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.0abd8208a0eb44c32d9f52395849381cb7361d5f.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/684e0c5928f94a0012b5261340bed1e0c0d7575f/e2e/6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md", "", "", "965104be-0ef5-4edb-82b1-facb0f37c968.md")
$zh.Hyperlinks.Add($zh.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/684e0c5928f94a0012b5261340bed1e0c0d7575f/e2e/6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/936ed7e448591055b0e50645239880596f4d0bd0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6b9ff258-5cf7-40c6-93dc-10784c1d9a46.0abd8208a0eb44c32d9f52395849381cb7361d5f.zh-cn.xlf", "", "", "965104be-0ef5-4edb-82b1-facb0f37c968.a05aeaf165ae5501d28d8aece37eedefb60075ad.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/684e0c5928f94a0012b5261340bed1e0c0d7575f/e2e/965104be-0ef5-4edb-82b1-facb0f37c968.md", "", "", "253d137b-9592-410f-9fca-d89327456d1f.md")
$zh.Hyperlinks.Add($zh.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/684e0c5928f94a0012b5261340bed1e0c0d7575f/e2e/965104be-0ef5-4edb-82b1-facb0f37c968.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/936ed7e448591055b0e50645239880596f4d0bd0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/965104be-0ef5-4edb-82b1-facb0f37c968.a05aeaf165ae5501d28d8aece37eedefb60075ad.zh-cn.xlf", "", "", "253d137b-9592-410f-9fca-d89327456d1f.fc74a366d75ceb48353102f324d1f380db212c73.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A3").Value() = "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md"
$de.Range("B3").Value() = ".md"
$de.Range("C3").Value() = "In Translation"
$de.Range("D3").Value() = "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.0abd8208a0eb44c32d9f52395849381cb7361d5f.de-de.xlf"
$de.Range("E3").Value() = "2016-03-11 16:25:17"
$de.Range("H3").Value() = "0001-01-01 00:00:00"
$de.Range("I3").Value() = "Include"

$de.Range("A4").Value() = "965104be-0ef5-4edb-82b1-facb0f37c968.md"
$de.Range("B4").Value() = ".md"
$de.Range("C4").Value() = "In Translation"
$de.Range("D4").Value() = "965104be-0ef5-4edb-82b1-facb0f37c968.a05aeaf165ae5501d28d8aece37eedefb60075ad.de-de.xlf"
$de.Range("E4").Value() = "2016-03-11 16:25:17"
$de.Range("H4").Value() = "0001-01-01 00:00:00"
$de.Range("I4").Value() = "Include"

$de.Range("A5").Value() = "253d137b-9592-410f-9fca-d89327456d1f.md"
$de.Range("B5").Value() = ".md"
$de.Range("C5").Value() = "Ready for handoff"
$de.Range("D5").Value() = "253d137b-9592-410f-9fca-d89327456d1f.fc74a366d75ceb48353102f324d1f380db212c73.de-de.xlf"
$de.Range("E5").Value() = "2016-03-11 16:24:14"
$de.Range("H5").Value() = "0001-01-01 00:00:00"
$de.Range("I5").Value() = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/99cf95aca13517e6b2d75943be5e9e0527602912/e2e/ef4783de-cfb1-4ba8-ade8-09b399ca2752.md", "", "", "ef4783de-cfb1-4ba8-ade8-09b399ca2752.md")
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/99cf95aca13517e6b2d75943be5e9e0527602912/e2e/ef4783de-cfb1-4ba8-ade8-09b399ca2752.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/48c6fa2a696c393ba6fe0406b8e403feafe97812/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ef4783de-cfb1-4ba8-ade8-09b399ca2752.69d912e4e3936c4a3dc4afbd5ae619737a6dcfd6.de-de.xlf", "", "", "ef4783de-cfb1-4ba8-ade8-09b399ca2752.69d912e4e3936c4a3dc4afbd5ae619737a6dcfd6.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8794364d964b3196e191cea2a456731fc66eb42a/e2e/ef4783de-cfb1-4ba8-ade8-09b399ca2752.md", "", "", "ef4783de-cfb1-4ba8-ade8-09b399ca2752.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/437ed5a90b07f1897d448046fc3cc2241be163b9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ef4783de-cfb1-4ba8-ade8-09b399ca2752.69d912e4e3936c4a3dc4afbd5ae619737a6dcfd6.de-de.xlf", "", "", "ef4783de-cfb1-4ba8-ade8-09b399ca2752.69d912e4e3936c4a3dc4afbd5ae619737a6dcfd6.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4daf793447831388b2108a25df80716aaad753b6/e2e/253d137b-9592-410f-9fca-d89327456d1f.md", "", "", "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md")
$de.Hyperlinks.Add($de.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/4daf793447831388b2108a25df80716aaad753b6/e2e/253d137b-9592-410f-9fca-d89327456d1f.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e328e07d1d5610e0ec79e93dab15b6fbd640fd09/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/253d137b-9592-410f-9fca-d89327456d1f.fc74a366d75ceb48353102f324d1f380db212c73.de-de.xlf", "", "", "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.0abd8208a0eb44c32d9f52395849381cb7361d5f.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/684e0c5928f94a0012b5261340bed1e0c0d7575f/e2e/6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md", "", "", "965104be-0ef5-4edb-82b1-facb0f37c968.md")
$de.Hyperlinks.Add($de.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/684e0c5928f94a0012b5261340bed1e0c0d7575f/e2e/6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6f8dcdc2dd2312f8cbc130df6559f4acc04640e3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6b9ff258-5cf7-40c6-93dc-10784c1d9a46.0abd8208a0eb44c32d9f52395849381cb7361d5f.de-de.xlf", "", "", "965104be-0ef5-4edb-82b1-facb0f37c968.a05aeaf165ae5501d28d8aece37eedefb60075ad.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/684e0c5928f94a0012b5261340bed1e0c0d7575f/e2e/965104be-0ef5-4edb-82b1-facb0f37c968.md", "", "", "253d137b-9592-410f-9fca-d89327456d1f.md")
$de.Hyperlinks.Add($de.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/684e0c5928f94a0012b5261340bed1e0c0d7575f/e2e/965104be-0ef5-4edb-82b1-facb0f37c968.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6f8dcdc2dd2312f8cbc130df6559f4acc04640e3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/965104be-0ef5-4edb-82b1-facb0f37c968.a05aeaf165ae5501d28d8aece37eedefb60075ad.de-de.xlf", "", "", "253d137b-9592-410f-9fca-d89327456d1f.fc74a366d75ceb48353102f324d1f380db212c73.de-de.xlf")
